$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.128.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.595.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "191.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("E7").Value = "  -1.89%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.591.60"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.181"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.666"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "56.02"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.73%  "
$ws.Range("E13").Value = "  +5.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.172.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.587.61"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.102.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.56%  "
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "476.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.11%  "
$ws.Range("E23").Value = "  +10.17%  "
$ws.Range("E24").Value = "  -6.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "95.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.18%  "
$ws.Range("E33").Value = "  +1.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "66.55"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "586.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "39.17"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.81%  "
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0805"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.396"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +18.78%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.47"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.24%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.138"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.228.72"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.28%  "
$ws.Range("E44").Value = "  +6.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0446"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.34%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.998"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.59%  "
